$wb = $excel.ActiveWorkbook

# --- Work on the "binek" sheet ---
$ws1 = $wb.Worksheets.Item("binek")

# Update existing values
$ws1.Range("B4").Value = 0.13
$ws1.Range("B5").Value = 0.08

# Add new row 7
$ws1.Range("A7").Value = "kredi_talep_esnekligi"
$ws1.Range("B7").Value = 0.63

# Select/activate the binek sheet and a specific range, mirroring the
# recorded selection state in the saved file.
$ws1.Select()
$ws1.Range("E19").Select()

$wb.Save()
